$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2561.1924
$ws.Range("I62").Value = 2553.4285
$ws.Range("J62").Value = 2570.25
$ws.Range("K62").Value = 2553.4285
$ws.Range("L62").Value = 2570.25
$ws.Range("M62").Value = -1929.4285
$ws.Range("N62").Value = -3818.25
$ws.Range("H64").Value = 3547.2632
$ws.Range("I64").Value = 3307.946
$ws.Range("K64").Value = 3307.946
$ws.Range("M64").Value = -3059.946
$ws.Range("H65").Value = 2561.1924
$ws.Range("I65").Value = 2553.4285
$ws.Range("J65").Value = 2570.25
$ws.Range("K65").Value = 12767.1425
$ws.Range("L65").Value = 12851.25
$ws.Range("M65").Value = -9647.1425
$ws.Range("N65").Value = -19091.25
$ws.Range("H67").Value = 3547.2632
$ws.Range("I67").Value = 3307.946
$ws.Range("K67").Value = 3307.946
$ws.Range("M67").Value = -2449.946
$ws.Range("H98").Value = 546.5909
$ws.Range("I98").Value = 546.5909
$ws.Range("K98").Value = 546.5909
$ws.Range("M98").Value = 951.4091
$ws.Range("H111").Value = 35715310
$ws.Range("I111").Value = 90910240
$ws.Range("J111").Value = 945.5294
$ws.Range("K111").Value = 272730720
$ws.Range("L111").Value = 2836.5882
$ws.Range("M111").Value = -272727653
$ws.Range("N111").Value = -8970.5882
$ws.Range("H116").Value = 3348172.8
$ws.Range("I116").Value = 3666718
$ws.Range("J116").Value = 3450
$ws.Range("K116").Value = 3666718
$ws.Range("L116").Value = 3450
$ws.Range("M116").Value = -3663276
$ws.Range("N116").Value = -10334
$ws.Range("H122").Value = 546.5909
$ws.Range("I122").Value = 546.5909
$ws.Range("K122").Value = 1639.7727
$ws.Range("M122").Value = 810.2273
$ws.Range("H137").Value = 40665.5
$ws.Range("I137").Value = 1682.8889
$ws.Range("J137").Value = 128376.375
$ws.Range("K137").Value = 5048.6667
$ws.Range("L137").Value = 385129.125
$ws.Range("M137").Value = -2498.6667
$ws.Range("N137").Value = -390229.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 34000
$ws.Range("J52").Value = 34000
$ws.Range("L52").Value = 34000
$ws.Range("N52").Value = -34526
$ws.Range("H121").Value = 34000
$ws.Range("J121").Value = 34000
$ws.Range("L121").Value = 34000
$ws.Range("N121").Value = -37494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1316.9375
$ws.Range("I31").Value = 1353.8572
$ws.Range("J31").Value = 1288.2222
$ws.Range("K31").Value = 1353.8572
$ws.Range("L31").Value = 1288.2222
$ws.Range("M31").Value = -1058.8572
$ws.Range("N31").Value = -1878.2222
$ws.Range("H34").Value = 1316.9375
$ws.Range("I34").Value = 1353.8572
$ws.Range("J34").Value = 1288.2222
$ws.Range("K34").Value = 1353.8572
$ws.Range("L34").Value = 1288.2222
$ws.Range("M34").Value = -1151.8572
$ws.Range("N34").Value = -1692.2222
$ws.Range("H99").Value = 2095.1155
$ws.Range("I99").Value = 1665.5555
$ws.Range("J99").Value = 3061.625
$ws.Range("K99").Value = 1665.5555
$ws.Range("L99").Value = 3061.625
$ws.Range("M99").Value = -167.5554999999999
$ws.Range("N99").Value = -6057.625
$ws.Range("H126").Value = 2095.1155
$ws.Range("I126").Value = 1665.5555
$ws.Range("J126").Value = 3061.625
$ws.Range("K126").Value = 4996.666499999999
$ws.Range("L126").Value = 9184.875
$ws.Range("M126").Value = -2526.666499999999
$ws.Range("N126").Value = -14124.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1295.8
$ws.Range("J68").Value = 1268.4688
$ws.Range("L68").Value = 3805.4064
$ws.Range("N68").Value = -5427.4064
$ws.Range("H71").Value = 1295.8
$ws.Range("J71").Value = 1268.4688
$ws.Range("L71").Value = 11416.2192
$ws.Range("N71").Value = -19528.2192
$ws.Range("H92").Value = 1250560.2
$ws.Range("I92").Value = 3333600.8
$ws.Range("J92").Value = 736
$ws.Range("K92").Value = 10000802.4
$ws.Range("L92").Value = 2208
$ws.Range("M92").Value = -9999554.399999999
$ws.Range("N92").Value = -4704
$ws.Range("H103").Value = 350
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 350
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 1050
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -2808
$ws.Range("H131").Value = 898.2857
$ws.Range("I131").Value = 305.55554
$ws.Range("J131").Value = 1059.9395
$ws.Range("K131").Value = 916.66662
$ws.Range("L131").Value = 3179.8185
$ws.Range("M131").Value = 4123.33338
$ws.Range("N131").Value = -13259.8185
$ws.Range("H136").Value = 2761.5
$ws.Range("I136").Value = 1774.6666
$ws.Range("J136").Value = 4876.143
$ws.Range("K136").Value = 5323.9998
$ws.Range("L136").Value = 14628.429
$ws.Range("M136").Value = -223.9997999999996
$ws.Range("N136").Value = -24828.429
$ws.Range("H138").Value = 1392.7273
$ws.Range("I138").Value = 968
$ws.Range("J138").Value = 2302.8572
$ws.Range("K138").Value = 2904
$ws.Range("L138").Value = 6908.571599999999
$ws.Range("M138").Value = 2236
$ws.Range("N138").Value = -17188.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 263.7
$ws.Range("I107").Value = 324.66666
$ws.Range("J107").Value = 237.57143
$ws.Range("K107").Value = 324.66666
$ws.Range("L107").Value = 237.57143
$ws.Range("M107").Value = 1595.33334
$ws.Range("N107").Value = -4077.57143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 501088.25
$ws.Range("I22").Value = 667740.3
$ws.Range("J22").Value = 1132
$ws.Range("K22").Value = 667740.3
$ws.Range("L22").Value = 1132
$ws.Range("M22").Value = -667445.3
$ws.Range("N22").Value = -1722
$ws.Range("H27").Value = 501088.25
$ws.Range("I27").Value = 667740.3
$ws.Range("J27").Value = 1132
$ws.Range("K27").Value = 667740.3
$ws.Range("L27").Value = 1132
$ws.Range("M27").Value = -667633.3
$ws.Range("N27").Value = -1346
$ws.Range("H40").Value = 4350.727
$ws.Range("I40").Value = 1738.5
$ws.Range("J40").Value = 5843.4287
$ws.Range("K40").Value = 1738.5
$ws.Range("L40").Value = 5843.4287
$ws.Range("M40").Value = -1602.5
$ws.Range("N40").Value = -6115.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8008.6665
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8008.6665
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 8008.6665
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -9880.666499999999
$ws.Range("H77").Value = 8008.6665
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8008.6665
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 24025.9995
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -33385.99950000001
$ws.Range("H126").Value = 2549.375
$ws.Range("I126").Value = 3597.5
$ws.Range("J126").Value = 1501.25
$ws.Range("K126").Value = 10792.5
$ws.Range("L126").Value = 4503.75
$ws.Range("M126").Value = -8322.5
$ws.Range("N126").Value = -9443.75
